# Apply crypto price/volume updates per diff (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.183.35"
$ws.Range("E2").Value = "  +0.68%  "

$ws.Range("D3").Value = "'3.521.63"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'593.52"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("D6").Value = "'174.08"
$ws.Range("E6").Value = "  +2.88%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.595"
$ws.Range("E8").Value = "  +3.79%  "

$ws.Range("E9").Value = "  +7.24%  "

$ws.Range("D10").Value = "'7.26"
$ws.Range("E10").Value = "  -0.47%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").Value = "'4.128.47"
$ws.Range("E12").Value = "  +0.42%  "

$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("D14").Value = "'29.13"
$ws.Range("E14").Value = "  +3.36%  "

$ws.Range("E15").Value = "  +1.54%  "

$ws.Range("D16").Value = "'67.179.38"
$ws.Range("E16").Value = "  +0.72%  "

$ws.Range("D17").Value = "'3.511.27"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").Value = "'14.26"
$ws.Range("E19").Value = "  +1.37%  "

$ws.Range("D20").Value = "'395.48"
$ws.Range("E20").Value = "  +1.21%  "

$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").Value = "'73.04"
$ws.Range("E22").Value = "  -0.53%  "

$ws.Range("D23").Value = "'0.542"
$ws.Range("E23").Value = "  +1.53%  "

$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("E25").Value = "  -3.50%  "

$ws.Range("D26").Value = "'10.28"
$ws.Range("E26").Value = "  +1.54%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").Value = "'6.28"
$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("E30").Value = "  -1.11%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("E34").Value = "  +1.76%  "

$ws.Range("D35").Value = "'163.17"
$ws.Range("E35").Value = "  +1.20%  "

$ws.Range("E36").Value = "  +0.85%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").Value = "'6.91"
$ws.Range("E38").Value = "  +3.72%  "

$ws.Range("D39").Value = "'4.71"
$ws.Range("E39").Value = "  +1.17%  "

$ws.Range("D40").Value = "'27.62"
$ws.Range("E40").Value = "  +4.43%  "

$ws.Range("D41").Value = "'0.0745"
$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").Value = "'26.44"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("D43").Value = "'2.65"
$ws.Range("E43").Value = "  +4.02%  "

$ws.Range("D44").Value = "'2.807.00"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("D45").Value = "'42.96"
$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").Value = "'337.13"
$ws.Range("E47").Value = "  -4.80%  "

$ws.Range("E48").Value = "  +1.17%  "

$ws.Range("D49").Value = "'33.54"
$ws.Range("E49").Value = "  -1.97%  "

$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("D51").Value = "'0.850"
$ws.Range("E51").Value = "  -0.48%  "

